# Cahier de bord: "TP de base :Faire marcher les 2 capteurs :"
#                -> "TP de base : Faire marcher les 2 capteurs :"
#
# The original paragraph is 3 runs wrapping a (now stale) grammar-check
# mark from Word's proofer:
#   [TP de base] <gramStart/> [ :Faire] <gramEnd/> [ marcher les 2 capteurs :]
# The fix simply adds the missing space after "base :" before "Faire", and
# the target markup drops the proofErr bookmarks and re-splits the text
# into three fresh runs:
#   [TP de base :] [ ] [Faire marcher les 2 capteurs :]
#
# Note: the colons here are preceded by a non-breaking space (U+00A0,
# French typographic spacing), not a regular space - preserve that.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

$oldText = "TP de base" + $nbsp + ":Faire marcher les 2 capteurs" + $nbsp + ":"
$newText = "TP de base" + $nbsp + ": Faire marcher les 2 capteurs" + $nbsp + ":"

# Replace the whole sentence in one shot: Word's Find/Replace rewrites the
# matched span as plain text, which naturally clears out the stale
# <w:proofErr/> gramStart/gramEnd marks that bracket the old run split.
$find = $d.Content
$find.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# Re-locate the sentence now that it reads correctly, and split off the
# newly-inserted space into its own run (matching the target run layout)
# by toggling a character formatting property across just that one
# character - Word always breaks a run at a formatting boundary, even a
# transient one that nets out to the same look.
$locate = $d.Content
$locate.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentenceStart = $locate.Start

$insertedSpace = $d.Range($sentenceStart + 12, $sentenceStart + 13)
$insertedSpace.Bold = $true
$insertedSpace.Bold = $false
